$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook already contains four "MY PINT ... v1.0" rows (246-249) that
# use the "peppol-doctype-wildcard" identifier scheme. This change adds four
# new, parallel rows (250-253) for the same MY PINT profiles but using the
# "busdox-docid-qns" identifier scheme instead, referencing a new comment
# (TICC-328).
#
# Because rows 246-249 already share almost every column value with the rows
# we need to add (Profile name, Doc type id value, Initial release, State,
# Issued by OpenPeppol?, BIS version, Domain Community, Associated process),
# the simplest and most faithful way to reproduce the exact formatting is to
# copy each of those rows down to the new row, then only touch the two
# columns that actually differ (B = identifier scheme, H = comment).

$ws.Rows.Item(246).Copy()
$ws.Rows.Item(250).Insert()

$ws.Rows.Item(247).Copy()
$ws.Rows.Item(251).Insert()

$ws.Rows.Item(248).Copy()
$ws.Rows.Item(252).Insert()

$ws.Rows.Item(249).Copy()
$ws.Rows.Item(253).Insert()

# Rows 247/249 (and therefore the copies in 251/253) render their long URN
# text wrapped across two lines, so restore that taller row height.
$ws.Rows.Item(251).RowHeight = 30
$ws.Rows.Item(253).RowHeight = 30

# Column B switches scheme from "peppol-doctype-wildcard" to
# "busdox-docid-qns" for the new rows; pick up the left-aligned formatting
# already used elsewhere in the sheet for that scheme (e.g. B238) before
# writing the new value.
$ws.Cells.Item(238, 2).Copy()
$ws.Range("B250:B253").PasteSpecial(-4122)
$ws.Range("B250:B253").Value = "busdox-docid-qns"

# Column H (Comment) references the new TICC ticket for this change.
$ws.Range("H250:H253").Value = "TICC-328"

# Restore the autofilter over the header + original data rows (unchanged
# range), since it otherwise gets dropped on save.
[void]$ws.Range("A1:L245").AutoFilter(1)

$excel.CutCopyMode = 0
